$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.368.22"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "2.281.59"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "157.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15,639.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "95.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "35.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0803"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").Value = "2.634.61"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").Value = "2.270.05"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.802"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.91%  "
$ws.Range("D19").Value = "42.257.86"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("D21").Value = "0.0₃0919"
$ws.Range("E21").Value = "  +1.97%  "
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "242.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0750"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("E37").Value = "  +5.21%  "
$ws.Range("E38").Value = "  +3.67%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  +6.21%  "
$ws.Range("D43").Value = "2.001.53"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  +12.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("E48").Value = "  +5.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
